$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert two new worksheets right after "Study":
#      "Stack1.Temperature"
#      "Virtual Stack1.Temperature CSV"
# ---------------------------------------------------------------------
$study = $wb.Worksheets.Item("Study")

$stackTemp = $wb.Worksheets.Add($null, $study)
$stackTemp.Name = "Stack1.Temperature"

$virtualStackCsv = $wb.Worksheets.Add($null, $stackTemp)
$virtualStackCsv.Name = "Virtual Stack1.Temperature CSV"

# ---------------------------------------------------------------------
# 2. Populate "Stack1.Temperature"
# ---------------------------------------------------------------------
$stackTemp.Range("A1").Value = "Identifier"
$stackTemp.Range("B1").Value = "Stack1.Temperature"
$stackTemp.Range("A2").Value = "VirtualInstrument"
$stackTemp.Range("B2").Value = "Virtual Stack1.Temperature CSV"
$stackTemp.Range("A4").Value = "Attributes.DataLogger.Destination"
$stackTemp.Range("A5").Value = "Attributes.DataLogger.Class"
$stackTemp.Range("A8").Value = "Attributes.Alarming.Limits.Hi-Hi"
$stackTemp.Range("B8").Value = 50
$stackTemp.Range("A9").Value = "Attributes.Alarming.Limits.Hi"
$stackTemp.Range("B9").Value = 45
$stackTemp.Range("A10").Value = "Attributes.Alarming.Limits.Lo"
$stackTemp.Range("B10").Value = 0
$stackTemp.Range("A11").Value = "Attributes.Alarming.Limits.Lo-Lo"
$stackTemp.Range("B11").Value = -1

$stackTemp.Range("A3").Select()

# ---------------------------------------------------------------------
# 3. Populate "Virtual Stack1.Temperature CSV"
# ---------------------------------------------------------------------
$virtualStackCsv.Range("A1").Value = "Identifier"
$virtualStackCsv.Range("B1").Value = "Virtual Stack1.Thermometer CSV"
$virtualStackCsv.Range("A2").Value = "Class"
$virtualStackCsv.Range("B2").Value = "Virtual Sensor from CSV.lvclass"
$virtualStackCsv.Range("A3").Value = "ConfigureClass"
$virtualStackCsv.Range("A5").Value = "Attributes.CsvParserClass"
$virtualStackCsv.Range("A6").Value = "Attributes.CsvFilepath"
$virtualStackCsv.Range("B6").Value = "D:\git\HAL\HAL-Devices\Tests\resources\2024-11-27-17-48-50_Imponator_Last Hope.csv"
$virtualStackCsv.Range("A7").Value = "Attributes.Input.Timestamp.Field"
$virtualStackCsv.Range("B7").Value = "Time in sec"
$virtualStackCsv.Range("A8").Value = "Attributes.Input.Timestamp.Format"
$virtualStackCsv.Range("B8").Value = "%t"
$virtualStackCsv.Range("A9").Value = "Attributes.Input.Value.Field"
$virtualStackCsv.Range("B9").Value = "Stack1.Temperature"

$virtualStackCsv.Range("B7").Select()

# ---------------------------------------------------------------------
# 4. Update "Study" sheet rows 4-18 (columns B and C)
# ---------------------------------------------------------------------
$study.Range("B4").Value = "Virtual Stack1.Temperature CSV"
$study.Range("B5").Value = "#Virtual Stack2.Temperature CSV"
$study.Range("B6").Value = "#Virtual Stack3.Temperature CSV"
$study.Range("B7").Value = "#Virtual Stack4.Temperature CSV"
$study.Range("B8").Value = "#Virtual Stack5.Temperature CSV"
$study.Range("B9").Value = "#Virtual Stack6.Temperature CSV"
$study.Range("B10").Value = "#Virtual Stack1.Fan CSV"
$study.Range("B11").Value = "#Virtual Stack2.Fan CSV"
$study.Range("B12").Value = "#Virtual Stack3.Fan CSV"
$study.Range("B13").Value = "#Virtual Stack4.Fan CSV"
$study.Range("B14").Value = "#Virtual Stack5.Fan CSV"
$study.Range("B15").Value = "#Virtual Stack6.Fan CSV"

$study.Range("C7").Value = "#Stack1.Temperature"
$study.Range("C8").Value = "#Stack2.Temperature"
$study.Range("C9").Value = "#Stack3.Temperature"
$study.Range("C10").Value = "#Stack4.Temperature"
$study.Range("C11").Value = "#Stack5.Temperature"
$study.Range("C12").Value = "#Stack6.Temperature"
$study.Range("C13").Value = "#Stack1.Fan"
$study.Range("C14").Value = "#Stack2.Fan"
$study.Range("C15").Value = "#Stack3.Fan"
$study.Range("C16").Value = "#Stack4.Fan"
$study.Range("C17").Value = "#Stack5.Fan"
$study.Range("C18").Value = "#Stack6.Fan"

$study.Range("B16").Select()

# ---------------------------------------------------------------------
# 5. Defined names local to "Stack1.Temperature" (localSheetId = 1),
#    mirroring the pattern used by the other MainSupply.* sheets.
# ---------------------------------------------------------------------
$stackTemp.Names.Add("Identifier", "=Stack1.Temperature!`$B`$1")
$stackTemp.Names.Add("VirtualInstrument", "=Stack1.Temperature!`$B`$2")
$stackTemp.Names.Add("Attributes.DataLogger.Destination", "=Stack1.Temperature!`$B`$4")
$stackTemp.Names.Add("Attributes.DataLogger.Class", "=Stack1.Temperature!`$B`$5")
$stackTemp.Names.Add("Attributes.Alarming.Limits.HiHi", "=Stack1.Temperature!`$B`$8")
$stackTemp.Names.Add("Attributes.Alarming.Limits.Hi", "=Stack1.Temperature!`$B`$9")
$stackTemp.Names.Add("Attributes.Alarming.Limits.Lo", "=Stack1.Temperature!`$B`$10")
$stackTemp.Names.Add("Attributes.Alarming.Limits.LoLo", "=Stack1.Temperature!`$B`$11")

# ---------------------------------------------------------------------
# 6. Active tab -> "Virtual Stack1.Temperature CSV" (index 2, 0-based)
# ---------------------------------------------------------------------
$virtualStackCsv.Activate()
